# PSDR2_BOM.xlsx update — add 8 new BOM rows (Amp Board parts) just above the
# "Number of different parts" / "Part count" summary rows, and move the
# selection/scroll position further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Push the two summary rows (old 89:91) down by inserting 8 blank rows above
# them. Excel copies the formatting of the row immediately above the insert
# point (row 88, style s="10" in column B) onto every new row, and the
# COUNT/SUM formulas in the (now) 97:98 rows keep referencing B2:B85 since
# that range is untouched by the insert.
$ws.Rows("89:96").Insert()

# --- Row 89: "Amp Board" header/part line (entered J, H, C like the source) ---
$ws.Range("J89").Value = "Amp Board"
$ws.Range("H89").Value = "LM3410XMF/NOPBCT-ND"
$ws.Range("C89").Value = "IC DRVR WT/OLED BCKLT SOT23-5"
$ws.Range("B89").Value = 1
$ws.Rows("89").RowHeight = 30

# --- Row 90: 0805 inductor, 10uH ---
$ws.Range("F90").Value = "10uH"
$ws.Range("C90").Value = "0805 Inductors"
$ws.Range("H90").Value = "535-10520-1-ND"
$ws.Range("B90").Value = 0

# --- Row 91: 0805 inductor, 470nH ---
$ws.Range("F91").Value = "470nH"
$ws.Range("C91").Value = "0805 Inductors"
$ws.Range("H91").Value = "535-10509-1-ND"
$ws.Range("B91").Value = 0

# --- Row 92: 2.2uH inductor (description cell keeps the default/no style) ---
$ws.Range("C92").Value = "INDUCTOR 2.2UH 260MA 20% SMD"
$ws.Range("C92").Style = "Normal"
$ws.Range("F92").Value = "2.2uH"
$ws.Range("H92").Value = "587-2043-1-ND"
$ws.Range("B92").Value = 2

# --- Row 93: 1.74 ohm resistor ---
$ws.Range("C93").Value = "RES 1.74 OHM 1/10W 1% 0603 SMD"
$ws.Range("H93").Value = "541-1.74HHCT-ND"
$ws.Range("B93").Value = 2
$ws.Rows("93").RowHeight = 30

# --- Row 94: 15uH inductor ---
$ws.Range("C94").Value = "INDUCTOR 15UH 20% 0805 SMD"
$ws.Range("H94").Value = "587-3020-1-ND"
$ws.Range("B94").Value = 1

# --- Row 95: Schottky diode ---
$ws.Range("C95").Value = "DIODE SCHOTTKY 40V SOD323"
$ws.Range("H95").Value = "SD101CWSTPMSCT-ND"
$ws.Range("B95").Value = 1

# Row 96 is left blank (only the inherited B96 formatting), matching source.

# Move the active selection the way the author last left it.
$ws.Range("I35").Select()
